# Regenerate merged AHB files
#
# The sheet contains a "diff" table comparing an old (FV2404) message
# implementation guide against a new (FV2410) one. The column headers
# used to be generically named "*_old" / "*_new"; they are renamed to
# carry the explicit format-version tags "*_FV2404" / "*_FV2410".
# The data range is then turned into a proper Excel Table, and the
# header row is frozen so it stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the header strings:
#      "<Name>_old" -> "<Name>_FV2404"
#      "<Name>_new" -> "<Name>_FV2410"
#    (the "diff" header in column K is left untouched)
[void]$ws.Cells.Replace("_old", "_FV2404")
[void]$ws.Cells.Replace("_new", "_FV2410")

# 2) Turn the used range A1:U63 into an Excel Table (ListObject) so the
#    header row drives an AutoFilter and the columns get named ranges.
$rng = $ws.Range("A1:U63")
$tbl = $ws.ListObjects.Add(1, $rng, 0, 1)
$tbl.Name = "Table1"

# 3) Freeze the header row (row 1) so it remains visible on scroll.
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
